$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Add the "Notes" sheet right after Sheet1, with the explanatory note in A1.
#    (This must happen before the B55 edit below so the shared-string table
#    ends up in the same order as the target workbook.)
$notes = $wb.Worksheets.Add($null, $ws1)
$notes.Name = "Notes"
$notes.Range("A1").Value = "IMPACT conversion for poultry removed because the USDA edible share correction is for removing bones."

# 2. Remove the IMPACT_conversion value for cpoul (row 9, poultry) -- the
#    note above explains why: the USDA edible-share correction already
#    accounts for bone removal, so a separate conversion factor is wrong.
$ws1.Range("C9").Clear()

# 3. Add the IMPACT_conversion value for cteas (row 40, tea).
#    Build the exact font/style combination (size-11 Calibri, theme color,
#    General number format) the same way Excel would: create it on a
#    scratch cell outside any styled column, then copy only the formatting
#    onto C40 so the column's "0.0" number format isn't re-applied.
$ws1.Range("Z1").Font.Size = 11
$ws1.Range("Z1").Value = "x"
$ws1.Range("C40").Value = 66.67
$ws1.Range("Z1").Copy()
$ws1.Range("C40").PasteSpecial(-4122)
$ws1.Range("Z1").Clear()

# 4. ctols (row 55): drop the "02024" entry from the USDA-code list.
$ws1.Range("B55").Value = "12104, 09193, 12021, 12023, 02033, 12160, 12012"

# 5. Restore the selection/active sheet to Sheet1 at B55.
$ws1.Activate()
$ws1.Range("B55").Select()
